$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.214.23"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "1.905.26"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'328.07"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "'0.4642"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("D8").Value = "'0.3956"
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("D9").Value = "'46.77"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").Value = "'0.07962"
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("D11").Value = "'0.9996"
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("D12").Value = "'22.28"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "1.924.48"
$ws.Range("E13").Value = "  +3.82%  "
$ws.Range("D14").Value = "'7.133"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").Value = "'5.767"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "'0.06972"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "'88.71"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "'1.006"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "'0.00001011"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").Value = "'17.15"
$ws.Range("E20").Value = "  +1.97%  "
$ws.Range("D21").Value = "'1.004"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "29.244.54"
$ws.Range("E22").Value = "  +2.12%  "
$ws.Range("D23").Value = "'5.354"
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("D24").Value = "'11.07"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "2.156.42"
$ws.Range("E25").Value = "  +3.63%  "
$ws.Range("D26").Value = "'2.050"
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("D27").Value = "'156.97"
$ws.Range("E27").Value = "  +2.98%  "
$ws.Range("D28").Value = "'19.50"
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("D29").Value = "'5.900"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").Value = "'119.27"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "'0.09412"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").Value = "'0.9230"
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("D34").Value = "'5.348"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").Value = "'1.347"
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("D37").Value = "'0.05829"
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("D38").Value = "'1.170"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").Value = "'0.02105"
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("D40").Value = "'8.000"
$ws.Range("E40").Value = "  +4.08%  "
$ws.Range("D41").Value = "'0.5753"
$ws.Range("E41").Value = "  +2.29%  "
$ws.Range("D42").Value = "'0.1811"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("D43").Value = "'9.979"
$ws.Range("E43").Value = "  +2.10%  "
$ws.Range("D44").Value = "'12.05"
$ws.Range("E44").Value = "  +2.72%  "
$ws.Range("D45").Value = "'0.5428"
$ws.Range("E45").Value = "  +2.66%  "
$ws.Range("D46").Value = "'2.218"
$ws.Range("E46").Value = "  +3.16%  "
$ws.Range("D47").Value = "'0.07101"
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("D48").Value = "'1.877"
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("D49").Value = "'2.598"
$ws.Range("E49").Value = "  +7.91%  "
$ws.Range("D50").Value = "'112.02"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").Value = "'1.064"
$ws.Range("E51").Value = "  -5.61%  "
